# Update the cached "datetimeFigureOut" Date placeholder text from
# 4/17/2022 to 4/18/2022 across the Slide Master, every Slide Layout,
# and the Notes Master.

$p = $ppt.ActivePresentation

$oldDate = "4/17/2022"
$newDate = "4/18/2022"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Slide Layout hanging off the master
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes
}

# Notes Master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes
